$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 8, pushing existing rows 8-42 down to 10-44.
$ws.Range("A8:A9").EntireRow.Insert()

# Populate new row 8
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = 44459
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 100112021
$ws.Range("G8").Value = "Ají"
$ws.Range("H8").Value = "Cristal"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 28000
$ws.Range("L8").Value = 30000
$ws.Range("M8").Value = 29000
$ws.Range("N8").Value = "`$/caja 15 kilos"
$ws.Range("O8").Value = "Región de Arica y Parinacota"
$ws.Range("P8").Value = 1933
$ws.Range("Q8").Value = 15
$ws.Range("R8").Value = "Hortaliza"

# Populate new row 9
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44459
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112021
$ws.Range("G9").Value = "Ají"
$ws.Range("H9").Value = "Inferno"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 130
$ws.Range("K9").Value = 38000
$ws.Range("L9").Value = 40000
$ws.Range("M9").Value = 39000
$ws.Range("N9").Value = "`$/caja 15 kilos"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 2600
$ws.Range("Q9").Value = 15
$ws.Range("R9").Value = "Hortaliza"
